# Generate Report for Handback
# The "bdc28621-516e-4ab6-92f9-1cb06dee1b24" file has now been handed back
# in sync with en-US, so update its status/datetime/error cells across all
# three report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-24 04:46:26"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-24 04:46:33"
$wsDeDe.Range("P3").Value = ""
